# Applies the NATMI TPM re-run for Gdf11-Acvr2a (xl/worksheets/sheet1.xml + xl/sharedStrings.xml):
#  - A new "Inflammatory-Mac" sending/target cluster is inserted (between FAPs and MuSCs)
#  - All Ligand/Receptor/Edge expression statistics are refreshed for the new TPM values
#  - The data block grows from 20 rows (A2:T21) to 25 rows (A2:T26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One array entry per data row (columns A..T, in order). Row 2 is the first entry, etc.
$data = @(
    @("ECs", "Gdf11", "Acvr2a", "ECs", 3, 1, 1.314013333333333, 3.94204, 0.2008810138079581, 0.2289960686502567, 3, 1, 9.24193, 27.72579, 0.1468938537243544, 0.1569651396557324, 12.14401924573333, 109.2961732116, 0.02950818625830621, 0.03594439989630123),
    @("ECs", "Gdf11", "Acvr2a", "FAPs", 3, 1, 1.314013333333333, 3.94204, 0.2008810138079581, 0.2289960686502567, 3, 1, 29.54200233333333, 88.626007, 0.469548954544906, 0.5017420086455576, 38.81858495936444, 349.36726463428, 0.09432347002144754, 0.1148969474565158),
    @("ECs", "Gdf11", "Acvr2a", "Inflammatory-Mac", 3, 1, 1.314013333333333, 3.94204, 0.2008810138079581, 0.2289960686502567, 3, 1, 7.349831333333333, 22.049494, 0.1168203014713749, 0.1248296948454213, 9.65777636975111, 86.91998732775998, 0.02346698059292108, 0.02858550937041269),
    @("ECs", "Gdf11", "Acvr2a", "MuSCs", 3, 1, 1.314013333333333, 3.94204, 0.2008810138079581, 0.2289960686502567, 2, 1, 12.1104985, 24.220997, 0.1924876941491673, 0.1371233128688515, 15.91335650231333, 95.48013901387998, 0.03866712314624088, 0.03140069956726615),
    @("ECs", "Gdf11", "Acvr2a", "Resolving-Mac", 3, 1, 1.314013333333333, 3.94204, 0.2008810138079581, 0.2289960686502567, 3, 1, 4.671440333333334, 14.014321, 0.07424919611019735, 0.079339843984437, 6.138334883871112, 55.24501395484, 0.01491525378904234, 0.01816851235976079),
    @("FAPs", "Gdf11", "Acvr2a", "ECs", 3, 1, 2.338584333333333, 7.015753, 0.3575132609679819, 0.4075503687484767, 3, 1, 9.24193, 27.72579, 0.1468938537243544, 0.1569651396557324, 21.61303270776333, 194.51729436987, 0.05251650066114768, 0.06397120054734987),
    @("FAPs", "Gdf11", "Acvr2a", "FAPs", 3, 1, 2.338584333333333, 7.015753, 0.3575132609679819, 0.4075503687484767, 3, 1, 29.54200233333333, 88.626007, 0.469548954544906, 0.5017420086455576, 69.08646383203011, 621.778174488271, 0.167869977923456, 0.2044851406400984),
    @("FAPs", "Gdf11", "Acvr2a", "Inflammatory-Mac", 3, 1, 2.338584333333333, 7.015753, 0.3575132609679819, 0.4075503687484767, 3, 1, 7.349831333333333, 22.049494, 0.1168203014713749, 0.1248296948454213, 17.18820040877578, 154.693803678982, 0.04176480692629395, 0.05087438816501127),
    @("FAPs", "Gdf11", "Acvr2a", "MuSCs", 3, 1, 2.338584333333333, 7.015753, 0.3575132609679819, 0.4075503687484767, 2, 1, 12.1104985, 24.220997, 0.1924876941491673, 0.1371233128688515, 28.32142206095683, 169.928532365741, 0.06881690323147631, 0.05588465672371317),
    @("FAPs", "Gdf11", "Acvr2a", "Resolving-Mac", 3, 1, 2.338584333333333, 7.015753, 0.3575132609679819, 0.4075503687484767, 3, 1, 4.671440333333334, 14.014321, 0.07424919611019735, 0.079339843984437, 10.92455717763478, 98.32101459871301, 0.02654507222560785, 0.03233498267230391),
    @("Inflammatory-Mac", "Gdf11", "Acvr2a", "ECs", 2, 0.6666666666666666, 0.3516466666666667, 1.05494, 0.05375831211924975, 0.06128225808512899, 3, 1, 9.24193, 27.72579, 0.1468938537243544, 0.1569651396557324, 3.249893878066667, 29.2490449026, 0.007896765636913263, 0.009619178198750907),
    @("Inflammatory-Mac", "Gdf11", "Acvr2a", "FAPs", 2, 0.6666666666666666, 0.3516466666666667, 1.05494, 0.05375831211924975, 0.06128225808512899, 3, 1, 29.54200233333333, 88.626007, 0.469548954544906, 0.5017420086455576, 10.38834664717555, 93.49511982458, 0.02524215925369247, 0.03074788326596808),
    @("Inflammatory-Mac", "Gdf11", "Acvr2a", "Inflammatory-Mac", 2, 0.6666666666666666, 0.3516466666666667, 1.05494, 0.05375831211924975, 0.06128225808512899, 3, 1, 7.349831333333333, 22.049494, 0.1168203014713749, 0.1248296948454213, 2.584543688928889, 23.26089320036, 0.006280062228363021, 0.007649845576205004),
    @("Inflammatory-Mac", "Gdf11", "Acvr2a", "MuSCs", 2, 0.6666666666666666, 0.3516466666666667, 1.05494, 0.05375831211924975, 0.06128225808512899, 2, 1, 12.1104985, 24.220997, 0.1924876941491673, 0.1371233128688515, 4.258616429196666, 25.55169857518, 0.01034781354118562, 0.008403226248716847),
    @("Inflammatory-Mac", "Gdf11", "Acvr2a", "Resolving-Mac", 2, 0.6666666666666666, 0.3516466666666667, 1.05494, 0.05375831211924975, 0.06128225808512899, 3, 1, 4.671440333333334, 14.014321, 0.07424919611019735, 0.079339843984437, 1.642696421748889, 14.78426779574, 0.003991511459095374, 0.004862124795488137),
    @("MuSCs", "Gdf11", "Acvr2a", "ECs", 2, 1, 2.409312, 4.818624, 0.3683258189716586, 0.2799174925428902, 3, 1, 9.24193, 27.72579, 0.1468938537243544, 0.1569651396557324, 22.26669285216, 133.60015711296, 0.05410479897492587, 0.04393728830907719),
    @("MuSCs", "Gdf11", "Acvr2a", "FAPs", 2, 1, 2.409312, 4.818624, 0.3683258189716586, 0.2799174925428902, 3, 1, 29.54200233333333, 88.626007, 0.469548954544906, 0.5017420086455576, 71.175900725728, 427.055404354368, 0.1729470032300386, 0.1404463649634976),
    @("MuSCs", "Gdf11", "Acvr2a", "Inflammatory-Mac", 2, 1, 2.409312, 4.818624, 0.3683258189716586, 0.2799174925428902, 3, 1, 7.349831333333333, 22.049494, 0.1168203014713749, 0.1248296948454213, 17.708036829376, 106.248220976256, 0.0430279332119602, 0.03494201517602447),
    @("MuSCs", "Gdf11", "Acvr2a", "MuSCs", 2, 1, 2.409312, 4.818624, 0.3683258189716586, 0.2799174925428902, 2, 1, 12.1104985, 24.220997, 0.1924876941491673, 0.1371233128688515, 29.17796936203199, 116.711877448128, 0.07089818758945818, 0.03838321390742314),
    @("MuSCs", "Gdf11", "Acvr2a", "Resolving-Mac", 2, 1, 2.409312, 4.818624, 0.3683258189716586, 0.2799174925428902, 3, 1, 4.671440333333334, 14.014321, 0.07424919611019735, 0.079339843984437, 11.254957252384, 67.529743514304, 0.02734789596527573, 0.02220861018686771),
    @("Resolving-Mac", "Gdf11", "Acvr2a", "ECs", 1, 0.3333333333333333, 0.1276956666666667, 0.383087, 0.01952159413315168, 0.02225381197324759, 3, 1, 9.24193, 27.72579, 0.1468938537243544, 0.1569651396557324, 1.180154412636667, 10.62138971373, 0.002867602193061398, 0.003493072704253217),
    @("Resolving-Mac", "Gdf11", "Acvr2a", "FAPs", 1, 0.3333333333333333, 0.1276956666666667, 0.383087, 0.01952159413315168, 0.02225381197324759, 3, 1, 29.54200233333333, 88.626007, 0.469548954544906, 0.5017420086455576, 3.772385682623223, 33.951471143609, 0.009166344116271341, 0.01116567231947781),
    @("Resolving-Mac", "Gdf11", "Acvr2a", "Inflammatory-Mac", 1, 0.3333333333333333, 0.1276956666666667, 0.383087, 0.01952159413315168, 0.02225381197324759, 3, 1, 7.349831333333333, 22.049494, 0.1168203014713749, 0.1248296948454213, 0.9385416119975557, 8.446874507978, 0.002280518511836602, 0.002777936557767879),
    @("Resolving-Mac", "Gdf11", "Acvr2a", "MuSCs", 1, 0.3333333333333333, 0.1276956666666667, 0.383087, 0.01952159413315168, 0.02225381197324759, 2, 1, 12.1104985, 24.220997, 0.1924876941491673, 0.1371233128688515, 1.546458179623167, 9.278749077738999, 0.003757666640806278, 0.003051516421732223),
    @("Resolving-Mac", "Gdf11", "Acvr2a", "Resolving-Mac", 1, 0.3333333333333333, 0.1276956666666667, 0.383087, 0.01952159413315168, 0.02225381197324759, 3, 1, 4.671440333333334, 14.014321, 0.07424919611019735, 0.079339843984437, 0.5965226876585557, 5.368704188927, 0.001449462671176057, 0.00176561397001646)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($v in $row) {
        $ws.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
    $r = $r + 1
}